$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number (e.g. "0.999") need to be
# pre-formatted as Text so Excel stores them as strings (matching the source
# data, which is all text) instead of auto-converting to a numeric cell.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "61.732.21"
$ws.Range("E2").Value = "  -1.20%  "
$ws.Range("D3").Value = "2.994.12"
$ws.Range("E3").Value = "  -2.23%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "531.95"
$ws.Range("E5").Value = "  -0.49%  "
$ws.Range("D6").Value = "132.06"
$ws.Range("E6").Value = "  -0.25%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "2.988.78"
$ws.Range("E8").Value = "  -2.19%  "
$ws.Range("D9").Value = "0.487"
$ws.Range("E9").Value = "  +0.18%  "
$ws.Range("D10").Value = "0.150"
$ws.Range("E10").Value = "  -1.56%  "
$ws.Range("D11").Value = "6.15"
$ws.Range("E11").Value = "  +1.03%  "
$ws.Range("D12").Value = "0.442"
$ws.Range("E12").Value = "  -2.24%  "
$ws.Range("D13").Value = "0.0000218"
$ws.Range("E13").Value = "  -1.85%  "
$ws.Range("D14").Value = "33.41"
$ws.Range("E14").Value = "  -2.06%  "
$ws.Range("D15").Value = "3.461.63"
$ws.Range("E15").Value = "  -1.33%  "
$ws.Range("E16").Value = "  +0.10%  "
$ws.Range("D17").Value = "61.613.23"
$ws.Range("E17").Value = "  -1.19%  "
$ws.Range("D18").Value = "2.978.62"
$ws.Range("E18").Value = "  -2.64%  "
$ws.Range("D19").Value = "6.51"
$ws.Range("E19").Value = "  -0.39%  "
$ws.Range("D20").Value = "459.04"
$ws.Range("E20").Value = "  -3.22%  "
$ws.Range("D21").Value = "13.16"
$ws.Range("E21").Value = "  -0.35%  "
$ws.Range("D22").Value = "0.674"
$ws.Range("E22").Value = "  -2.80%  "
$ws.Range("D23").Value = "6.84"
$ws.Range("E23").Value = "  -3.51%  "
$ws.Range("D24").Value = "77.64"
$ws.Range("E24").Value = "  -0.89%  "
$ws.Range("D25").Value = "11.83"
$ws.Range("E25").Value = "  -0.79%  "
$ws.Range("D26").Value = "0.997"
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").Value = "2.65"
$ws.Range("E27").Value = "  -0.85%  "
$ws.Range("D28").Value = "7.66"
$ws.Range("E28").Value = "  -6.27%  "
$ws.Range("D29").Value = "0.997"
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("D30").Value = "25.52"
$ws.Range("E30").Value = "  -0.24%  "
$ws.Range("E31").Value = "  +3.38%  "
$ws.Range("D32").Value = "1.83"
$ws.Range("E32").Value = "  -0.50%  "
$ws.Range("D33").Value = "55.54"
$ws.Range("E33").Value = "  -1.69%  "
$ws.Range("D34").Value = "5.38"
$ws.Range("E34").Value = "  +4.09%  "
$ws.Range("D35").Value = "2.24"
$ws.Range("E35").Value = "  -5.64%  "
$ws.Range("D36").Value = "5.79"
$ws.Range("E36").Value = "  -1.48%  "
$ws.Range("D37").Value = "453.70"
$ws.Range("E37").Value = "  -3.22%  "
$ws.Range("D38").Value = "3.135.88"
$ws.Range("E38").Value = "  +1.06%  "
$ws.Range("D39").Value = "0.0384"
$ws.Range("E39").Value = "  -0.84%  "
$ws.Range("D40").Value = "0.0779"
$ws.Range("E40").Value = "  -0.84%  "
$ws.Range("D41").Value = "0.117"
$ws.Range("E41").Value = "  +4.58%  "
$ws.Range("D42").Value = "7.98"
$ws.Range("E42").Value = "  -0.05%  "
$ws.Range("D43").Value = "2.41"
$ws.Range("E43").Value = "  -6.12%  "
$ws.Range("D45").Value = "0.244"
$ws.Range("E45").Value = "  -0.97%  "
$ws.Range("D46").Value = "25.33"
$ws.Range("E46").Value = "  +5.35%  "
$ws.Range("D47").Value = "120.74"
$ws.Range("E47").Value = "  +3.32%  "
$ws.Range("D48").Value = "0.107"
$ws.Range("E48").Value = "  +0.84%  "
$ws.Range("D49").Value = "1.93"
$ws.Range("E49").Value = "  -3.92%  "
$ws.Range("D50").Value = "0.0₃0505"
$ws.Range("E50").Value = "  -0.99%  "
$ws.Range("D51").Value = "1.23"
$ws.Range("E51").Value = "  +6.34%  "

# Restore the default (General) cell style now that the text value is locked in,
# so the only observable change is the cell content, not its formatting.
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D51").Style = "Normal"

